# Apply weekly update to the "Ají" sheet: insert two new price records at row 21,
# pushing the existing historical rows (old rows 21-132) down to rows 23-134.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 21 (shifts everything below down by 2 rows)
$ws.Rows("21:22").Insert()

# Common/constant values for this market (same across every data row in the sheet)
$mercadoId   = 1
$mercado     = "Agrícola del Norte S.A. de Arica"
$region      = "Arica y Parinacota"
$codreg      = 15
$categoriaId = 100112021
$categoria   = "Ají"
$unidad      = "$/caja 15 kilos"
$origen      = "Región de Arica y Parinacota"
$kgUnidades  = 15
$clasif      = "Hortaliza"

# New row 21: Ají Inferno, Primera
$r = 21
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = 45030
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $categoriaId
$ws.Cells.Item($r, 7).Value2  = $categoria
$ws.Cells.Item($r, 8).Value2  = "Inferno"
$ws.Cells.Item($r, 9).Value2  = "Primera"
$ws.Cells.Item($r, 10).Value2 = 160
$ws.Cells.Item($r, 11).Value2 = 17000
$ws.Cells.Item($r, 12).Value2 = 18000
$ws.Cells.Item($r, 13).Value2 = 17500
$ws.Cells.Item($r, 14).Value2 = $unidad
$ws.Cells.Item($r, 15).Value2 = $origen
$ws.Cells.Item($r, 16).Value2 = 1167
$ws.Cells.Item($r, 17).Value2 = $kgUnidades
$ws.Cells.Item($r, 18).Value2 = $clasif

# New row 22: Ají Inferno, Segunda
$r = 22
$ws.Cells.Item($r, 1).Value2  = $mercadoId
$ws.Cells.Item($r, 2).Value2  = $mercado
$ws.Cells.Item($r, 3).Value2  = $region
$ws.Cells.Item($r, 4).Value2  = 45030
$ws.Cells.Item($r, 5).Value2  = $codreg
$ws.Cells.Item($r, 6).Value2  = $categoriaId
$ws.Cells.Item($r, 7).Value2  = $categoria
$ws.Cells.Item($r, 8).Value2  = "Inferno"
$ws.Cells.Item($r, 9).Value2  = "Segunda"
$ws.Cells.Item($r, 10).Value2 = 130
$ws.Cells.Item($r, 11).Value2 = 14000
$ws.Cells.Item($r, 12).Value2 = 15000
$ws.Cells.Item($r, 13).Value2 = 14500
$ws.Cells.Item($r, 14).Value2 = $unidad
$ws.Cells.Item($r, 15).Value2 = $origen
$ws.Cells.Item($r, 16).Value2 = 967
$ws.Cells.Item($r, 17).Value2 = $kgUnidades
$ws.Cells.Item($r, 18).Value2 = $clasif
